$wb = $excel.ActiveWorkbook

# Delete worksheet named "123"
$ws3 = $wb.Worksheets.Item("123")
$excel.DisplayAlerts = $false
$ws3.Delete()

# Update header cell text in the two remaining sheets
$ws1 = $wb.Worksheets.Item("121")
$ws1.Range("C3").Value = "N_ZACHET"

$ws2 = $wb.Worksheets.Item("122")
$ws2.Range("C3").Value = "N_ZACHET"
